$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text in E8 ("Good Morning" -> "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection recorded in the sheet view
$ws.Range("E8").Select()
